$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change column C values from numeric 4 to the text "space" for every 5th
# row starting at row 6 through row 301 (rows 6,11,16,...,301).
for ($r = 6; $r -le 301; $r += 5) {
    $ws.Range("C$r").Value = "space"
}

# Update the sheet's active selection / view position to F14 (also resets
# the scrolled-to top-left cell back to the default).
$ws.Range("F14").Select()
